$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")

# Updated Sudoku candidate/solution grid (A1:I9)
$gridValues = @{
    "A1" = 1
    "B1" = 6
    "C1" = 8
    "D1" = 3
    "E1" = 5
    "F1" = 2
    "G1" = 4
    "H1" = 9
    "I1" = 7
    "A2" = 9
    "B2" = 35
    "C2" = 2
    "D2" = 1
    "E2" = 7
    "F2" = 4
    "A3" = 35
    "B3" = 7
    "C3" = 4
    "D3" = 9
    "E3" = 8
    "F3" = 6
    "I3" = 1
    "A4" = 2
    "C4" = 69
    "D4" = 4
    "F4" = 3
    "G4" = 7
    "B5" = 1
    "C5" = 7
    "D5" = 8
    "E5" = 9
    "F5" = 5
    "A6" = 56
    "C6" = 3
    "D6" = 2
    "F6" = 7
    "I6" = 4
    "A7" = 8
    "B7" = 4
    "C7" = 5
    "D7" = 7
    "E7" = 2
    "F7" = 1
    "G7" = 369
    "H7" = 36
    "I7" = 39
    "A8" = 36
    "C8" = 69
    "D8" = 5
    "E8" = 4
    "F8" = 8
    "G8" = 1
    "H8" = 7
    "I8" = 2
    "A9" = 7
    "B9" = 2
    "C9" = 1
    "D9" = 6
    "E9" = 3
    "F9" = 9
    "G9" = 5
    "H9" = 4
    "I9" = 8
}
foreach ($addr in $gridValues.Keys) {
    $ws.Range($addr).Value = $gridValues[$addr]
}

# Cells that no longer hold a candidate list (cleared, style kept)
$clearedCells = @(
    "G2",
    "H2",
    "I2",
    "G3",
    "H3",
    "B4",
    "E4",
    "H4",
    "I4",
    "G5",
    "H5",
    "B6",
    "E6",
    "G6",
    "H6",
    "B8"
)
foreach ($addr in $clearedCells) {
    $ws.Range($addr).Value = $null
}

# Scratch note added below the grid
$ws.Range("E11").Value = 16

# Recalculate the COUNTIF-based completion formulas in column K
$excel.Calculate()

# Leave the selection where the author left off
$ws.Activate()
$ws.Range("E4").Select()
